# edit.ps1 - applies the MAUI EXAMPLE.docx README restructuring described
# by the commit "Updated to .NET 9.   Fixed some obsolete choices in the
# XAML code."
#
# Strategy: walk the document from the top, editing paragraphs in place
# where the target keeps the same logical slot, deleting paragraphs whose
# content is being relocated/rewritten, and inserting freshly-styled
# paragraphs for all of the new README sections, finishing with the
# trailing "Charles B Hayes" line-break tweak.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Paragraph 2 (blank paragraph right under the title) becomes a
#    "NoSpacing" styled blank paragraph instead of the default style.
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Style = "NoSpacing"

# ---------------------------------------------------------------------
# 2. Paragraph 3 (the long intro paragraph) becomes "NoSpacing" styled,
#    gets "VS 2022 " inserted right before ".NET Multi-platform App UI",
#    and loses its trailing run of spaces.
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Style = "NoSpacing"

$introRng = $p3.Range
$found = $introRng.Find.Execute(".NET Multi-platform App UI", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$introRng.Collapse(1)
$introRng.InsertBefore("VS 2022 ")

$p3 = $d.Paragraphs.Item(3)
$trimRng = $p3.Range
$trimRng.Find.Execute("trying to win.          ", $false, $false, $false, $false, $false, $true, 1, $false, "trying to win.", 1) | Out-Null

# ---------------------------------------------------------------------
# 3. Insert a new blank "NoSpacing" paragraph after the intro paragraph.
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3.Range.InsertParagraphAfter()
$pBlank1 = $d.Paragraphs.Item(4)
$pBlank1.Range.Style = "NoSpacing"

# ---------------------------------------------------------------------
# 4. Paragraph 5 ("Supports" heading, Heading1 style) is retitled to
#    "Possible Future Changes" - it keeps its Heading1 + bottom border
#    pPr untouched.
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Find.Execute("Supports", $false, $false, $false, $false, $false, $true, 1, $false, "Possible Future Changes", 1) | Out-Null

# ---------------------------------------------------------------------
# 5. Delete the old ".NET 7.0 for ..." list (paragraphs 6-9) plus the
#    trailing blank "NoSpacing" paragraph (10), the "License" heading
#    (11) and the "Open source - MIT License" paragraph (12). All of
#    this content gets rebuilt/relocated further down in the document.
# ---------------------------------------------------------------------
$pStart = $d.Paragraphs.Item(6)
$pEnd = $d.Paragraphs.Item(12)
$delRng = $d.Range($pStart.Range.Start, $pEnd.Range.End)
Write-Host "Deleting block:" $delRng.Text
$delRng.Delete()

Write-Host "Stage 2 paragraph count:" $d.Paragraphs.Count

# ---------------------------------------------------------------------
# 6. Build the new "Possible Future Changes" bullet-less list, the new
#    "Known Issues" section, and the relocated/renumbered "Supports"
#    (.NET 9 targets) and "License" sections. Everything is inserted
#    right after paragraph 5 ("Possible Future Changes"), in reverse
#    order, always anchored on paragraph 5 so indices never drift.
# ---------------------------------------------------------------------

function Insert-After($anchorIndex, $style, $text) {
    $anchor = $d.Paragraphs.Item($anchorIndex)
    $anchor.Range.InsertParagraphAfter()
    $newp = $d.Paragraphs.Item($anchorIndex + 1)
    # InsertParagraphAfter() copies the *direct* pPr formatting (e.g. the
    # Heading1 bottom pBdr) from the anchor paragraph, not just its style
    # reference. Always explicitly (re)apply a style - "Normal" clears
    # that inherited direct formatting for the style-less paragraphs.
    if ($style -ne $null) {
        $newp.Range.Style = $style
    } else {
        $newp.Range.Style = "Normal"
    }
    if ($text -ne $null) {
        $newp.Range.Text = $text
    }
    return $newp
}

# Insert in reverse order (bottom-most line first) so that paragraph 5
# stays the anchor index throughout and we don't have to recompute
# offsets after every insertion.

# "Open source - MIT License" (NoSpacing)
$enDash = [char]0x2013
Insert-After 5 "NoSpacing" ("Open source " + $enDash + " MIT License") | Out-Null
# blank NoSpacing before it
Insert-After 5 "NoSpacing" "" | Out-Null
# "License" Heading1
$pLicense = Insert-After 5 "Heading1" "License"
# blank NoSpacing before License heading
Insert-After 5 "NoSpacing" "" | Out-Null

# ".NET 9.0 for MacCatalyst (Untested)"
Insert-After 5 "NoSpacing" ".NET 9.0 for MacCatalyst (Untested)" | Out-Null
# ".NET 9.0 for IOS (Untested)"
Insert-After 5 "NoSpacing" ".NET 9.0 for IOS (Untested)" | Out-Null
# ".NET 9.0 for Android"
Insert-After 5 "NoSpacing" ".NET 9.0 for Android" | Out-Null
# ".NET 9.0 for Windows"
Insert-After 5 "NoSpacing" ".NET 9.0 for Windows" | Out-Null
# blank NoSpacing
Insert-After 5 "NoSpacing" "" | Out-Null
# "Supports" Heading1
$pSupports = Insert-After 5 "Heading1" "Supports"

# trailing-space-only paragraph (no style)
Insert-After 5 $null "          " | Out-Null
# "none"
Insert-After 5 "NoSpacing" "none" | Out-Null
# blank NoSpacing
Insert-After 5 "NoSpacing" "" | Out-Null
# "Known Issues" Heading1
$pKnownIssues = Insert-After 5 "Heading1" "Known Issues"
# blank paragraph, no style (plain <w:p/>)
Insert-After 5 $null "" | Out-Null

# "Add the ability for the computer to occasionally loose, currently it never does."
Insert-After 5 "NoSpacing" "Add the ability for the computer to occasionally loose, currently it never does." | Out-Null
# "Add a small delay before the computer plays."
Insert-After 5 "NoSpacing" "Add a small delay before the computer plays." | Out-Null
# "Add a shaking effect when the user clicks down on the X or O."
Insert-After 5 "NoSpacing" "Add a shaking effect when the user clicks down on the X or O." | Out-Null
# "Ability for the user to pick X or O."
Insert-After 5 "NoSpacing" "Ability for the user to pick X or O." | Out-Null
# "If I have time in the future, there are a few things that I was wanting to add:"
Insert-After 5 "NoSpacing" "If I have time in the future, there are a few things that I was wanting to add:" | Out-Null
# blank NoSpacing right after "Possible Future Changes"
Insert-After 5 "NoSpacing" "" | Out-Null

Write-Host "Stage 3 paragraph count:" $d.Paragraphs.Count
for ($i=1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs.Item($i)
  Write-Host $i ":" $p.Range.Text
}
